$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.874826908111572
$ws.Range("B1").Value = 4.147600173950195
$ws.Range("C1").Value = 3.638537168502808
$ws.Range("D1").Value = 4.537003040313721
$ws.Range("E1").Value = 5.081137180328369
